$wb = $excel.ActiveWorkbook

# --- Update selection on "branch4" sheet (D16) and drop it as the active tab ---
$branch4 = $wb.Worksheets.Item("branch4")
[void]$branch4.Range("D16").Select()

# --- Create the new sheet by duplicating "branch4" (keeps sheet formatting/namespaces consistent) then wipe its contents ---
$branch4.Copy($null, $branch4)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "branch4 original"
$ws.Cells.Clear()

# column width -> 65 characters
$ws.Columns.Item(1).ColumnWidth = 64.15

# --- Write cell values ---
$ws.Range("A1").Value = 'Test Cases'
$ws.Range("A2").Value = 'Verify if an existing user is able to login to Mars with valid email address and password'
$ws.Range("A3").Value = 'Verify if user is taken to their home page upon login in to Mars successfully with valid credentials '
$ws.Range("A4").Value = 'Verify if a new user is able to register to Mars using the register functionality'
$ws.Range("A5").Value = 'Verify if an existing user is not allowed to login to Mars with valid email address and invalid password'
$ws.Range("A6").Value = 'Verify if an existing user is not allowed to login to mars with invalid email address and valid password'
$ws.Range("A7").Value = 'Verify if an existing user is not allowed to login to Mars with invalid email address and invalid password'
$ws.Range("A8").Value = 'Verify if an existing user is not allowed to login to Mars with null email address and a valid password'
$ws.Range("A9").Value = 'Verify if an existing user is not allowed to login to Mars with valid email address and null password'
$ws.Range("A10").Value = 'Verify if an existing user is not allowed to login to Mars with null email address and null password'
$ws.Range("A11").Value = 'Verify if an existing user is not allowed to login to Mars with null email address and invalid password'
$ws.Range("A12").Value = 'Verify if an existing user is not allowed to login to Mars with invalid email address and null password'
$ws.Range("A13").Value = 'Verify if an existing user is not allowed to see other user''s home page upon login'
$ws.Range("A14").Value = 'Verify new user is able to add a language'
$ws.Range("A15").Value = 'Verify existing user is able to add a language '
$ws.Range("A16").Value = 'Verify new user is able to delete a language'
$ws.Range("A17").Value = 'Verify existing user is able to delete a language'
$ws.Range("A18").Value = 'Verify new user is able to update a language '

# --- Row heights (wrapped rows get an explicit 26pt height) ---
$ws.Rows.Item(2).RowHeight = 26
$ws.Rows.Item(3).RowHeight = 26
$ws.Rows.Item(5).RowHeight = 26
$ws.Rows.Item(6).RowHeight = 26
$ws.Rows.Item(7).RowHeight = 26
$ws.Rows.Item(8).RowHeight = 26
$ws.Rows.Item(9).RowHeight = 26
$ws.Rows.Item(10).RowHeight = 26
$ws.Rows.Item(11).RowHeight = 26
$ws.Rows.Item(12).RowHeight = 26
$ws.Rows.Item(13).RowHeight = 26

# --- Fonts: build the two custom fonts once, then reuse them ---
# Body font: Arial 10 (non-bold) -- established on A2 then reused everywhere else
$ws.Range("A2").Font.Name = "Arial"
$ws.Range("A2").Font.Size = 10
$ws.Range("A3").Font.Name = "Arial"
$ws.Range("A3").Font.Size = 10
$ws.Range("A4").Font.Name = "Arial"
$ws.Range("A4").Font.Size = 10
$ws.Range("A5").Font.Name = "Arial"
$ws.Range("A5").Font.Size = 10
$ws.Range("A6").Font.Name = "Arial"
$ws.Range("A6").Font.Size = 10
$ws.Range("A7").Font.Name = "Arial"
$ws.Range("A7").Font.Size = 10
$ws.Range("A8").Font.Name = "Arial"
$ws.Range("A8").Font.Size = 10
$ws.Range("A9").Font.Name = "Arial"
$ws.Range("A9").Font.Size = 10
$ws.Range("A10").Font.Name = "Arial"
$ws.Range("A10").Font.Size = 10
$ws.Range("A11").Font.Name = "Arial"
$ws.Range("A11").Font.Size = 10
$ws.Range("A12").Font.Name = "Arial"
$ws.Range("A12").Font.Size = 10
$ws.Range("A13").Font.Name = "Arial"
$ws.Range("A13").Font.Size = 10

# Header font: Arial 10 Bold
$ws.Range("A1").Font.Name = "Arial"
$ws.Range("A1").Font.Size = 10
$ws.Range("A1").Font.Bold = $true

# --- Fills ---
# Header row: solid gray fill (RGB 153,153,153 = FF999999)
$ws.Range("A1").Interior.Color = 10066329
# Rows 11-12: solid white fill (FFFFFFFF)
$ws.Range("A11").Interior.Color = 16777215
$ws.Range("A12").Interior.Color = 16777215

# --- Borders: thin box border around every used cell A1:A18 ---
$ws.Range("A1:A18").Borders.LineStyle = 1

# --- Wrap text on rows 1-13 (header + body groups) ---
$ws.Range("A1").WrapText = $true
$ws.Range("A2").WrapText = $true
$ws.Range("A3").WrapText = $true
$ws.Range("A4").WrapText = $true
$ws.Range("A5").WrapText = $true
$ws.Range("A6").WrapText = $true
$ws.Range("A7").WrapText = $true
$ws.Range("A8").WrapText = $true
$ws.Range("A9").WrapText = $true
$ws.Range("A10").WrapText = $true
$ws.Range("A11").WrapText = $true
$ws.Range("A12").WrapText = $true
$ws.Range("A13").WrapText = $true
# rows 14-18 explicitly have wrap text off (plain bordered cells)
$ws.Range("A14").WrapText = $false
$ws.Range("A15").WrapText = $false
$ws.Range("A16").WrapText = $false
$ws.Range("A17").WrapText = $false
$ws.Range("A18").WrapText = $false

[void]$ws.Range("D6").Select()

# --- Make the new sheet the active tab ---
$ws.Activate()

Write-Host "done"

